# Balance Report.xlsx - add "Цена"/"Сумма"/"Из них негодно" columns to the
# header row and resize the affected columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells. "Сумма" (G1) is written before "Из них негодно" (E1)
# so the shared-string table grows in the same order the source workbook
# shows (Сумма before Из них негодно).
$ws.Range("G1").Value = "Сумма"
$ws.Range("E1").Value = "Из них негодно"
$ws.Range("F1").Value = "Цена"
$ws.Range("H1").Value = "Из них негодно"

# Column widths: column A grew a bit, and the four new columns were sized
# to fit their headers.
$ws.Columns.Item(1).ColumnWidth = 8.2504
$ws.Columns.Item(5).ColumnWidth = 13.5834
$ws.Columns.Item(6).ColumnWidth = 11.2503
$ws.Columns.Item(7).ColumnWidth = 14.917
$ws.Columns.Item(8).ColumnWidth = 18.251

# Leave the cursor on F7, matching the saved selection.
$ws.Range("F7").Select() | Out-Null
